$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: spread (C2) and total (D2) inputs
$ws.Range("C2").Value = -6.5
$ws.Range("D2").Value = 34

# Row 3: total (D3) input
$ws.Range("D3").Value = 42

# Row 4: spread (C4) input
$ws.Range("C4").Value = -3

# Row 5: spread (C5) input
$ws.Range("C5").Value = -4
